$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-07-07")
$src.Copy([System.Reflection.Missing]::Value, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-07-08"

$data = New-Object "object[,]" 50,4
$data[0,0] = 1
$data[0,1] = "食い詰め傭兵の幻想奇譚"
$data[0,2] = "原作／まいん キャラクター原案／peroshi 漫画／池宮アレア"
$data[0,3] = "第28話"
$data[1,0] = 2
$data[1,1] = "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～"
$data[1,2] = "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり"
$data[1,3] = "第５０話　雌雄を決する器用貧乏（１）"
$data[2,0] = 3
$data[2,1] = "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！"
$data[2,2] = "島知宏 音速炒飯 有都あらゆる"
$data[2,3] = "第２２食　ユクシーさんの覚悟、すごいのですわ！（１）"
$data[3,0] = 4
$data[3,1] = "異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～"
$data[3,2] = "モリタ Ｕ４ nima"
$data[3,3] = "第12話（６）　恋焦がれ！黄金色の実りと、握るはふっくら銀のシャリ（６）"
$data[4,0] = 5
$data[4,1] = "濁る瞳で何を願う ハイセルク戦記"
$data[4,2] = "トルトネン 創-taro 斎藤八呑"
$data[4,3] = "第31話 大暴走"
$data[5,0] = 6
$data[5,1] = "ふかふかダンジョン攻略記～俺の異世界転生冒険譚～"
$data[5,2] = "KAKERU"
$data[5,3] = "第66話「東アイギス」（前半）"
$data[6,0] = 7
$data[6,1] = "ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ"
$data[6,2] = "漫画：晴野しゅー 原作：ちんくるり キャラクター原案：イセ川ヤスタカ"
$data[6,3] = "第71話前半"
$data[7,0] = 8
$data[7,1] = "いつでも自宅に帰れる俺は、異世界で行商人をはじめました"
$data[7,2] = "漫画／明地雫 原作／霜月緋色 キャラクター原案／いわさきたかし"
$data[7,3] = "第58話"
$data[8,0] = 9
$data[8,1] = "王都ワンオペゴーレムマスター。まさかの追放！？～自由の身になったので弟子の美人勇者たちと一緒に最強ゴーレム作ります。戻ってこいと言われてももう知らん！～@COMIC"
$data[8,2] = "阿住 周（漫画） レルクス（原作） 布施龍太（キャラクター原案）"
$data[8,3] = "第8話"
$data[9,0] = 10
$data[9,1] = "異世界のすみっこで快適ものづくり生活 ～女神さまのくれた工房はちょっとやりすぎ性能だった～"
$data[9,2] = "西山アラタ(漫画) 長田信織(原作) 東上文(キャラクター原案)"
$data[9,3] = "EP.19①"
$data[10,0] = 11
$data[10,1] = "ポーション頼みで生き延びます！ 続"
$data[10,2] = "原作：FUNA 漫画：園心ふつう キャラクター原案：すきま"
$data[10,3] = "第63話　長いお別れ Ⅱ"
$data[11,0] = 12
$data[11,1] = "江戸前エルフ"
$data[11,2] = "樋口彰彦"
$data[11,3] = "#112"
$data[12,0] = 13
$data[12,1] = "ブチ切れ令嬢は報復を誓いました。 ～魔導書の力で祖国を叩き潰します～"
$data[12,2] = "漫画：おおのいも 原作：はぐれメタボ キャラクター原案：昌未"
$data[12,3] = "第47話"
$data[13,0] = 14
$data[13,1] = "創造錬金術師は自由を謳歌する 故郷を追放されたら、魔王のお膝元で超絶効果のマジックアイテム作り放題になりました"
$data[13,2] = "姫乃タカ(漫画) 千月さかき(原作) かぼちゃ(キャラクター原案)"
$data[13,3] = "第19話-2"
$data[14,0] = 15
$data[14,1] = "クラス最安値で売られた俺は、実は最強パラメーター"
$data[14,2] = "カンブリア爆発太郎(漫画) RYOMA(原作) 黒井ススム(キャラクター原案)"
$data[14,3] = "第35話-3"
$data[15,0] = 16
$data[15,1] = "ロメリア戦記～伯爵令嬢、魔王を倒した後も人類やばそうだから軍隊組織する～"
$data[15,2] = "漫画：上戸 亮 原作：有山リョウ(小学館「ガガガブックス」刊) キャラクター原案：コダマ"
$data[15,3] = "第13話「助けてくれる人々」②"
$data[16,0] = 17
$data[16,1] = "孤児院テイマー"
$data[16,2] = "漫画：倉崎もろこ 原作：安藤正樹 キャラクター原案：イシバシヨウスケ"
$data[16,3] = "第60話"
$data[17,0] = 18
$data[17,1] = "断れない会長は友江くんにだけしてあげたい"
$data[17,2] = "沼地どろまる(著者)"
$data[17,3] = "コミックス第２巻発売記念！生徒会総選挙！"
$data[18,0] = 19
$data[18,1] = "アレクサンダー英雄戦記～最強の土魔術士～"
$data[18,2] = "マツオカヨシノリ るれくちぇ なんじゃもんじゃ"
$data[18,3] = "第10話（後編）"
$data[19,0] = 20
$data[19,1] = "願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜"
$data[19,2] = "ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)"
$data[19,3] = "第4話-1：師匠と弟子の新生活"
$data[20,0] = 21
$data[20,1] = "え、テイマーは使えないってパーティから追放したよね？ ～実は世界唯一の【精霊使い】だと判明した途端に手のひらを返されても遅い。精霊の王女様にめちゃくちゃ溺愛されながら、僕はマイペースに最強を目指すので"
$data[20,2] = "漫画：最中なつめ 原作：茨木野"
$data[20,3] = "第30話 特性と代償(後編)"
$data[21,0] = 22
$data[21,1] = "神の庭付き楠木邸"
$data[21,2] = "安斎アキラ(著者) えんじゅ(原作) ox(キャラクター原案)"
$data[21,3] = "第32話"
$data[22,0] = 23
$data[22,1] = "姫様“拷問”の時間です"
$data[22,2] = "原作:春原ロビンソン　漫画:ひらけい"
$data[22,3] = "拷問144"
$data[23,0] = 24
$data[23,1] = "Ｒｅ：ゼロから始める異世界生活 第四章 聖域と強欲の魔女"
$data[23,2] = "花鶏ハルノ(作画) 相川有(構成) 長月達平(原作) 大塚真一郎(キャラクター原案)"
$data[23,3] = "第61話①　エリオール大森林の永久凍土"
$data[24,0] = 25
$data[24,1] = "迷宮ブラックカンパニー"
$data[24,2] = "安村洋平"
$data[24,3] = "第50話　落花流水（前編）"
$data[25,0] = 26
$data[25,1] = "ジゼルの錬金飴"
$data[25,2] = "漫画： katoson 原作：斯波 キャラクター原案：LINO"
$data[25,3] = "第8話"
$data[26,0] = 27
$data[26,1] = "いとこのこ"
$data[26,2] = "いぬちく(著者)"
$data[26,3] = "第35.5話"
$data[27,0] = 28
$data[27,1] = "追放されたギルド職員は、世界最強の召喚士@COMIC"
$data[27,2] = "原作：月島秀一 漫画：あづち涼 キャラクター原案：チワワ丸"
$data[27,3] = "第9話②「伏魔殿ダラスの惨劇」"
$data[28,0] = 29
$data[28,1] = "コボルト無双、モフモフな最弱噛ませ犬だけど世界最強を目指す！"
$data[28,2] = "赤志木ひの乃 shiba"
$data[28,3] = "第十三話 帰還"
$data[29,0] = 30
$data[29,1] = "少年マールの転生冒険記"
$data[29,2] = "漫画家：あわや 原作：月ノ宮マクラ キャラクター原案：まっちょこ"
$data[29,3] = "第17話"
$data[30,0] = 31
$data[30,1] = "無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する"
$data[30,2] = "漫画：さらさみさ 小説： 桜霧琥珀 キャラ原案： 福きつね"
$data[30,3] = "第19話前半"
$data[31,0] = 32
$data[31,1] = "魔物ノ森ノ少女ノヴァ"
$data[31,2] = "kamatama"
$data[31,3] = "第6話「ロスタイム」後編"
$data[32,0] = 33
$data[32,1] = "このヒーラー、めんどくさい"
$data[32,2] = "丹念に発酵(著者)"
$data[32,3] = "第88話：ゴブリンの罠"
$data[33,0] = 34
$data[33,1] = "安達としまむら"
$data[33,2] = "柚原もけ(漫画) 入間人間(原作) のん(キャラクターデザイン)"
$data[33,3] = "第48話「最初の旅の端１」①"
$data[34,0] = 35
$data[34,1] = "スキル【万物支配】に目覚めたおっさんは、ダンジョンで生計を立てることにしました～無職から始める支配者無双～"
$data[34,2] = "岸本和葉 原田 臙 シミズヒロノリ 吉武"
$data[34,3] = "第3話　パーティ結成‼"
$data[35,0] = 36
$data[35,1] = "ポンコツスキルしか使えない悪役魔女だけど、テイムしたパリピなスライムたちと強く生きます！"
$data[35,2] = "漫画：鈴木イゾ 原作：雨傘ヒョウゴ キャラクター原案：朝日川日和"
$data[35,3] = "第8話"
$data[36,0] = 37
$data[36,1] = "ぽんドロイド！ はまさん"
$data[36,2] = "はれやまはれぞう(著者)"
$data[36,3] = "第3話"
$data[37,0] = 38
$data[37,1] = "ニャイト・オブ・ザ・リビングキャット"
$data[37,2] = "原作：ホークマン 作画：メカルーツ"
$data[37,3] = "Chapter15　ニャンペイジ（後編）"
$data[38,0] = 39
$data[38,1] = "異世界で最強のスキルを生み出せたので、ひたすら無双することにしました。　～俺だけがスキルの数値を勝手に操作～"
$data[38,2] = "漫画：星トマジロウ 原作：ヒゲ抜き地蔵 キャラクター原案：山椒魚"
$data[38,3] = "第9話 ②"
$data[39,0] = 40
$data[39,1] = "アラフォーおっさんはスローライフの夢を見るか？"
$data[39,2] = "漫画：大関詠詞 原作：サイトウアユム キャラクター原案： ジョンディー"
$data[39,3] = "第14話"
$data[40,0] = 41
$data[40,1] = "アイドル辞めるけど結婚してくれますか!?"
$data[40,2] = "三吉汐美(著者)"
$data[40,3] = "第16話前半"
$data[41,0] = 42
$data[41,1] = "まったく最近の探偵ときたら"
$data[41,2] = "五十嵐正邦(著者)"
$data[41,3] = "第113話"
$data[42,0] = 43
$data[42,1] = "勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが"
$data[42,2] = "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)"
$data[42,3] = "第3話 前編"
$data[43,0] = 44
$data[43,1] = "エンドロールの後には最高の旅を"
$data[43,2] = "リキタケ"
$data[43,3] = "最終話 最高の旅を②"
$data[44,0] = 45
$data[44,1] = "千年英雄"
$data[44,2] = "原作/福島航平 作画/中村ゆきひろ"
$data[44,3] = "18話①"
$data[45,0] = 46
$data[45,1] = "さんしょく弁当"
$data[45,2] = "兎月あい(著者)"
$data[45,3] = "第18話#3"
$data[46,0] = 47
$data[46,1] = "リアデイルの大地にて"
$data[46,2] = "月見だしお(著者) Ceez(原作) てんまそ(キャラクター原案) 涼風涼(構成)"
$data[46,3] = "第39章-2"
$data[47,0] = 48
$data[47,1] = "異世界創造のすゝめ～スマホアプリで惑星を創ってしまった俺は神となり世界を巡る～@COMIC"
$data[47,2] = "漫画：岩戸あきら 原作：たまごかけキャンディー キャラクター原案：かれい"
$data[47,3] = "第11話 ①"
$data[48,0] = 49
$data[48,1] = "王子様の友達"
$data[48,2] = "すけろく(著者)"
$data[48,3] = "第28話"
$data[49,0] = 50
$data[49,1] = "クソ女に幸あれ"
$data[49,2] = "岸川瑞樹"
$data[49,3] = "第59話"

$newSheet.Range("A2:D51").Value = $data
